# Update "想去人数" (want-to-go count) figures in the "展览" (Exhibitions)
# sheet and the "全部类型" (All types) sheet, reflecting the refreshed data
# output ("Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions only) -------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

$wsExpo.Range("F2").Value = 13431   # was 13385
$wsExpo.Range("F4").Value = 655     # was 650
$wsExpo.Range("F5").Value = 221     # was 220
$wsExpo.Range("F6").Value = 457     # was 452
$wsExpo.Range("F7").Value = 1344    # was 1330
$wsExpo.Range("F8").Value = 127     # was 125

# --- Sheet "全部类型" (all event types, includes the same exhibitions) ----
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Range("F2").Value = 13431    # was 13385
$wsAll.Range("F4").Value = 655      # was 650
$wsAll.Range("F5").Value = 221      # was 220
$wsAll.Range("F8").Value = 457      # was 452
$wsAll.Range("F9").Value = 1344     # was 1330
$wsAll.Range("F11").Value = 127     # was 125
